# The "2024" sheet tracks one SMS/notification entry per row in columns R
# (message/category text) and S (timestamp), ordered newest-first starting
# at row 44. A new, newer entry ("dispute" at 2024-09-19 22:27:16) arrived
# ahead of the existing top entry, so every existing entry from row 44
# downward shifts down by one row (which also carries the lone "Broadband"
# label in column A from row 172 down to row 173), and the new entry is
# written into the freshly inserted row 44.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at row 44, pushing rows 44:172 down to 45:173
# (dimension grows from A1:Y172 to A1:Y173 automatically).
$ws.Rows("44:44").Insert()

# Populate the newly inserted row with the latest entry.
$ws.Range("R44").Value = "dispute"
$ws.Range("S44").Value = "2024-09-19 22:27:16"
